# The deck's live DrawingML theme (ppt/theme/theme2.xml, the part wired to
# the slide master / presentation theme relationship) currently holds the
# "Integral" colour scheme. The edit swaps it back to the stock "Office
# Theme" colour scheme (the scheme that ships on ppt/theme/theme1.xml, the
# part used only by the notes master).
#
# PowerPoint's object model edits a presentation's active theme colours
# through Slide.ThemeColorScheme (or SlideRange.ThemeColorScheme) - the
# twelve slots are, in order: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
# ThemeColor.RGB uses the same packed 0xBBGGRR layout as VBA's RGB(r,g,b),
# so a plain 0xRRGGBB literal has to be byte-swapped before it is stored.

$p   = $ppt.ActivePresentation
$tcs = $p.Slides.Item(1).ThemeColorScheme

function HexToPpRGB([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r -bor ($g -shl 8) -bor ($b -shl 16)
}

# Office Theme colour scheme (target values, slot order matches
# ThemeColorScheme.Colors(1..12)).
$officeTheme = @(
    "000000",  #  1 dk1
    "FFFFFF",  #  2 lt1
    "44546A",  #  3 dk2
    "E7E6E6",  #  4 lt2
    "5B9BD5",  #  5 accent1
    "ED7D31",  #  6 accent2
    "A5A5A5",  #  7 accent3
    "FFC000",  #  8 accent4
    "4472C4",  #  9 accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

for ($i = 1; $i -le $tcs.Count; $i++) {
    $tcs.Colors($i).RGB = HexToPpRGB($officeTheme[$i - 1])
}
